# B6-PowerPoint.pptx edit: 6 Jun 2020
#
# 1) Re-style the three tables (slides 14, 15, 16) from the custom
#    "Table_0" style to the built-in "Medium Style 2 - Accent 1" table
#    style ({32301A4C-7D0D-41A4-A863-EB7BD86F3773}).
# 2) Swap the deck's colour theme from "Integral"/"Red Violet" to the
#    stock "Office Theme"/"Office" palette (font scheme + format scheme
#    are already identical between the two themes, so only the twelve
#    theme colours need to move).

$p = $ppt.ActivePresentation

# --- 1) Tables: apply the built-in table style -----------------------
$targetStyleId = "{32301A4C-7D0D-41A4-A863-EB7BD86F3773}"

for ($idx = 1; $idx -le $p.Slides.Count; $idx++) {
    $slide = $p.Slides.Item($idx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyleId)
        }
    }
}

# --- 2) Theme colours: switch to the "Office Theme" palette ----------
$themeColors = $p.SlideMaster.Theme.ThemeColorScheme

# index -> (theme element, target RGB hex "RRGGBB")
$officeThemeRgb = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

for ($i = 1; $i -le 12; $i++) {
    $hex = $officeThemeRgb[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # PowerPoint's RGB COM property is packed as 0x00BBGGRR.
    $bgr = ($b * 65536) + ($g * 256) + $r
    $themeColors.Colors($i).RGB = $bgr
}
